$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Confirm the features that are planned to implement by assigning "Colin"
# as the person assigned to the watermark-related tasks.
$ws.Range("C11").Value = "Colin"
$ws.Range("C12").Value = "Colin"

# Update the current selection to match the author's saved cursor position.
$ws.Range("C20").Select()
